# The block of 5 "Usage" columns (AE:AI) gets reordered. Each column's
# header label stays paired with its original value, but the columns
# move into a new left-to-right order:
#   before: tkm-N1Usage, tkm-N3Usage, tkm-N2Usage, pkmUsage, tkm-SZMUsage
#   after : tkm-N3Usage, pkmUsage,    tkm-SZMUsage, tkm-N2Usage, tkm-N1Usage

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("tkm-N3Usage", "pkmUsage", "tkm-SZMUsage", "tkm-N2Usage", "tkm-N1Usage")
$values  = @(130.3, 858, 414.5, 24.2, 7.5)

$startCol = 31  # column AE
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $values[$i]
}
